# Turns the blank "firsttable" workbook into a two-column header row
# ("First Name" / "Last name") styled like a MySQL-export table header,
# matching the target OOXML diff:
#   - xl/sharedStrings.xml gains "First Name" / "Last name"
#   - xl/worksheets/sheet1.xml gains the A1:B1 header row + A2 selection
#   - xl/styles.xml gains the two dxfs (header-row / whole-table) that back
#     the named "MySqlDefault" table style

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row cell values -------------------------------------------------
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last name"

# --- Register the table-style dxfs (bold/grey header, plain whole table) ---
# Using a throwaway conditional format is the only COM path that writes new
# <dxf> entries into xl/styles.xml; we delete the rule right after so the
# worksheet ends up with no <conditionalFormatting> left behind, same as the
# target workbook (the dxf registrations persist in the styles part).

# dxf #0 - header row: bold text on a light grey fill
$fcHeader = $ws.Range("A1:B1").FormatConditions.Add(1, 4, "1")
$fcHeader.Font.Bold = $true
$fcHeader.Interior.Color = 14145495
$fcHeader.Delete()

# dxf #1 - whole table: no special emphasis / no fill
$fcTable = $ws.Range("A1:B1").FormatConditions.Add(1, 4, "2")
$fcTable.Font.Bold = $false
$fcTable.Interior.Pattern = -4142
$fcTable.Delete()

# --- Sheet selection/dimension ---------------------------------------------
# Dimension grows to A1:B1 automatically from the writes above; leave the
# active selection on A2 like the committed file.
$ws.Range("A2").Select()
